$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared text updates (title block) ---
$ws.Range("A8").Value2 = "Volume 32   Number  48"
$ws.Range("C9").Value2 = "Report Covering the Week  11/24/2025  Through  11/30/2025"

# --- Style fixes: C15 and C27 change from text placeholder to numeric style (#,##0) ---
$ws.Range("C15").NumberFormat = "#,##0"
$ws.Range("C27").NumberFormat = "#,##0"

# --- Crime data numeric updates ---
$ws.Range("N14").Value2 = -71.428571428571
$ws.Range("C15").Value2 = 1
$ws.Range("D15").Value2 = 2
$ws.Range("E15").Value2 = -50
$ws.Range("G15").Value2 = 7
$ws.Range("H15").Value2 = -57.142857142857
$ws.Range("I15").Value2 = 28
$ws.Range("J15").Value2 = 38
$ws.Range("K15").Value2 = -26.315789473684
$ws.Range("L15").Value2 = 33.333333333333
$ws.Range("M15").Value2 = 40
$ws.Range("N15").Value2 = -17.647058823529
$ws.Range("C16").Value2 = 4
$ws.Range("D16").Value2 = 6
$ws.Range("E16").Value2 = -33.333333333333
$ws.Range("G16").Value2 = 22
$ws.Range("H16").Value2 = 18.181818181818
$ws.Range("I16").Value2 = 283
$ws.Range("J16").Value2 = 233
$ws.Range("K16").Value2 = 21.459227467811
$ws.Range("L16").Value2 = 10.116731517509
$ws.Range("M16").Value2 = -3.082191780821
$ws.Range("N16").Value2 = -54.792332268370
$ws.Range("C17").Value2 = 13
$ws.Range("D17").Value2 = 11
$ws.Range("E17").Value2 = 18.181818181818
$ws.Range("F17").Value2 = 35
$ws.Range("G17").Value2 = 34
$ws.Range("H17").Value2 = 2.941176470588
$ws.Range("I17").Value2 = 495
$ws.Range("J17").Value2 = 438
$ws.Range("K17").Value2 = 13.013698630137
$ws.Range("L17").Value2 = 28.571428571428
$ws.Range("M17").Value2 = 103.703703703704
$ws.Range("N17").Value2 = 54.6875
$ws.Range("C18").Value2 = 4
$ws.Range("D18").Value2 = 7
$ws.Range("E18").Value2 = -42.857142857142
$ws.Range("F18").Value2 = 16
$ws.Range("G18").Value2 = 14
$ws.Range("H18").Value2 = 14.285714285714
$ws.Range("I18").Value2 = 212
$ws.Range("J18").Value2 = 156
$ws.Range("K18").Value2 = 35.897435897435
$ws.Range("L18").Value2 = -2.304147465437
$ws.Range("M18").Value2 = -41.111111111111
$ws.Range("N18").Value2 = -84.525547445255
$ws.Range("C19").Value2 = 17
$ws.Range("D19").Value2 = 15
$ws.Range("E19").Value2 = 13.333333333333
$ws.Range("F19").Value2 = 57
$ws.Range("G19").Value2 = 53
$ws.Range("H19").Value2 = 7.547169811320
$ws.Range("I19").Value2 = 732
$ws.Range("J19").Value2 = 785
$ws.Range("K19").Value2 = -6.751592356687
$ws.Range("L19").Value2 = 16.006339144215
$ws.Range("M19").Value2 = 78.102189781021
$ws.Range("N19").Value2 = 35.055350553505
$ws.Range("C20").Value2 = 3
$ws.Range("D20").Value2 = 11
$ws.Range("E20").Value2 = -72.727272727272
$ws.Range("F20").Value2 = 28
$ws.Range("G20").Value2 = 38
$ws.Range("H20").Value2 = -26.315789473684
$ws.Range("I20").Value2 = 477
$ws.Range("J20").Value2 = 425
$ws.Range("K20").Value2 = 12.235294117647
$ws.Range("L20").Value2 = -1.851851851851
$ws.Range("M20").Value2 = 127.142857142857
$ws.Range("N20").Value2 = -73.005093378607
$ws.Range("C21").Value2 = 42
$ws.Range("D21").Value2 = 52
$ws.Range("E21").Value2 = -19.230769230769
$ws.Range("G21").Value2 = 168
$ws.Range("H21").Value2 = -1.785714285714
$ws.Range("I21").Value2 = 2231
$ws.Range("J21").Value2 = 2078
$ws.Range("K21").Value2 = 7.362848893166
$ws.Range("L21").Value2 = 11.438561438561
$ws.Range("M21").Value2 = 44.307891332470
$ws.Range("N21").Value2 = -52.257650331692
$ws.Range("F22").Value2 = 2
$ws.Range("G22").Value2 = 2
$ws.Range("I22").Value2 = 15
$ws.Range("J22").Value2 = 23
$ws.Range("K22").Value2 = -34.782608695652
$ws.Range("L22").Value2 = 36.363636363636
$ws.Range("M22").Value2 = -16.666666666666
$ws.Range("D23").Value2 = 3
$ws.Range("E23").Value2 = -66.666666666666
$ws.Range("F23").Value2 = 4
$ws.Range("G23").Value2 = 7
$ws.Range("H23").Value2 = -42.857142857142
$ws.Range("I23").Value2 = 102
$ws.Range("J23").Value2 = 99
$ws.Range("K23").Value2 = 3.030303030303
$ws.Range("L23").Value2 = -6.422018348623
$ws.Range("M23").Value2 = 72.881355932203
$ws.Range("C24").Value2 = 29
$ws.Range("E24").Value2 = 26.086956521739
$ws.Range("F24").Value2 = 166
$ws.Range("G24").Value2 = 90
$ws.Range("H24").Value2 = 84.444444444444
$ws.Range("I24").Value2 = 1665
$ws.Range("J24").Value2 = 1150
$ws.Range("K24").Value2 = 44.782608695652
$ws.Range("L24").Value2 = 13.265306122449
$ws.Range("M24").Value2 = 93.830034924330
$ws.Range("C25").Value2 = 12
$ws.Range("D25").Value2 = 7
$ws.Range("E25").Value2 = 71.428571428571
$ws.Range("F25").Value2 = 78
$ws.Range("G25").Value2 = 29
$ws.Range("H25").Value2 = 168.965517241379
$ws.Range("I25").Value2 = 690
$ws.Range("J25").Value2 = 435
$ws.Range("K25").Value2 = 58.620689655172
$ws.Range("L25").Value2 = 9.523809523809
$ws.Range("C26").Value2 = 11
$ws.Range("D26").Value2 = 17
$ws.Range("E26").Value2 = -35.294117647058
$ws.Range("F26").Value2 = 39
$ws.Range("G26").Value2 = 49
$ws.Range("H26").Value2 = -20.408163265306
$ws.Range("I26").Value2 = 594
$ws.Range("J26").Value2 = 530
$ws.Range("K26").Value2 = 12.075471698113
$ws.Range("L26").Value2 = 21.224489795918
$ws.Range("M26").Value2 = -0.834724540901
$ws.Range("C27").Value2 = 1
$ws.Range("D27").Value2 = 2
$ws.Range("E27").Value2 = -50
$ws.Range("G27").Value2 = 8
$ws.Range("H27").Value2 = -50
$ws.Range("I27").Value2 = 33
$ws.Range("J27").Value2 = 43
$ws.Range("K27").Value2 = -23.255813953488
$ws.Range("L27").Value2 = -2.941176470588
$ws.Range("D28").Value2 = 1
$ws.Range("E28").Value2 = 0
$ws.Range("F28").Value2 = 5
$ws.Range("G28").Value2 = 7
$ws.Range("H28").Value2 = -28.571428571428
$ws.Range("I28").Value2 = 66
$ws.Range("J28").Value2 = 64
$ws.Range("K28").Value2 = 3.125
$ws.Range("L28").Value2 = 4.761904761904
$ws.Range("M29").Value2 = -53.846153846153
$ws.Range("N29").Value2 = -76.470588235294
$ws.Range("M30").Value2 = -68.181818181818
$ws.Range("N30").Value2 = -84.090909090909
